$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Version 1.0" -> "Version 1.1"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Version 1.0", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Version 1.1", 2)

# ---------------------------------------------------------------------------
# 2. Drop "javascript " from the description paragraph and, per the target
#    revision, store the new sentence as three separate runs:
#       "Thi" | "s program is a simple" | " calculator program."
# ---------------------------------------------------------------------------
$oldText = "This program is a simple javascript calculator program."

$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $oldText) {
        $para = $p
        break
    }
}

$fullRange = $para.Range
# Exclude the trailing paragraph mark so only the run text is replaced,
# leaving the <w:p> element (and its rsid attributes) untouched.
$runsRange = $d.Range($fullRange.Start, $fullRange.End - 1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Thi</w:t></w:r><w:r><w:t>s program is a simple</w:t></w:r><w:r><w:t xml:space="preserve"> calculator program.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$runsRange.InsertXML($xml)
